# "Added sheet number argument to processFile"
# The existing col-oriented "person" sheet (Sheet1) is joined by a second
# person record placed on a new "Sheet2" so that processFile can be pointed
# at a specific sheet index.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Re-select Sheet1's full data range (matches the post-edit file: the
# previous edit had left the cursor parked on row 8) before the new sheet
# is added and becomes active.
$ws1.Range("A1:C12").Select()

# Insert the new sheet right after Sheet1 - it becomes "Sheet2" and the
# active/selected sheet, same as Worksheets.Add() does when Excel is asked
# to add a tab after an existing one.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Second person record, laid out the same column-oriented way as Sheet1
# (col A = field name, col B = value for this record).
$ws2.Range("A1").Value = "firstName"
$ws2.Range("B1").Value = "Max"

$ws2.Range("A2").Value = "lastName"
$ws2.Range("B2").Value = "Irwin"

$ws2.Range("A3").Value = "address.street"
$ws2.Range("B3").Value = "123 Fake Street"

$ws2.Range("A4").Value = "address.city"
$ws2.Range("B4").Value = "Rochester"

$ws2.Range("A5").Value = "address.state"
$ws2.Range("B5").Value = "NY"

$ws2.Range("A6").Value = "address.zip"
$ws2.Range("B6").Value = 99999

$ws2.Range("A7").Value = "isEmployee"
$ws2.Range("B7").Formula = '="false"'

$ws2.Range("A8").Value = "phones[0].type"
$ws2.Range("B8").Value = "home"

$ws2.Range("A9").Value = "phones[0].number"
$ws2.Range("B9").Value = "123.456.7890"

$ws2.Range("A10").Value = "phones[1].type"
$ws2.Range("B10").Value = "work"

$ws2.Range("A11").Value = "phones[1].number"
$ws2.Range("B11").Value = "505-505-1010"

$ws2.Range("A12").Value = "aliases[]"
$ws2.Range("B12").Value = "binarymax;arch"

# Column widths sized (auto-fit-like) to comfortably show the longest
# entries in each column.
$ws2.Columns.Item(1).ColumnWidth = 15.5
$ws2.Columns.Item(2).ColumnWidth = 16

# Put the cursor on the last populated cell, like a user who just finished
# typing the table would leave it.
$ws2.Range("B12").Select()

# Page setup matches Sheet2's (new-sheet) defaults.
$ws2.PageSetup.LeftMargin = 0.75 * 72
$ws2.PageSetup.RightMargin = 0.75 * 72
$ws2.PageSetup.TopMargin = 1 * 72
$ws2.PageSetup.BottomMargin = 1 * 72
$ws2.PageSetup.HeaderMargin = 0.5 * 72
$ws2.PageSetup.FooterMargin = 0.5 * 72

# Sheet1 picks up a portrait/letter-ish page setup too.
$ws1.PageSetup.PaperSize = 9
$ws1.PageSetup.Orientation = 1
